$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new column before EY (day "28-dec") ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Insert a new column, shifting EY:GC (and everything after) one column to the right.
$ws1.Range("EY1").EntireColumn.Insert()

# Row 1 header for the freshly inserted column.
$ws1.Range("EY1").Value = "28-dec"

# Data rows 2-25 get the same "-" placeholder used by the other empty days.
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 155).Value = "-"
}

# --- Sheet "Gaz": append the new last row ---
$ws2 = $wb.Worksheets.Item("Gaz")
# Leading apostrophe forces text storage (matches the other "yyyy-mm-dd" date
# strings in column A), avoiding Excel auto-converting it to a real date.
$ws2.Range("A183").Value = "'2025-12-26"
$ws2.Range("B183").Value = 27.75
